$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the timestamp recorded for the previous availability-check batch
#    (rows 324-337): 44232.15592659322 -> 44232.15592659722
# ---------------------------------------------------------------------------
for ($r = 324; $r -le 337; $r++) {
  $ws.Cells.Item($r, 4).Value = 44232.15592659722
}

# ---------------------------------------------------------------------------
# 2) Append a brand-new availability-check batch: rows 338-351, one row per
#    monitored service, all stamped with the same new check timestamp.
# ---------------------------------------------------------------------------
$newTimestamp = 44232.17696617641

$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
  "https://www.dataintelligence-group.com/",
  "https://serviciodashboard.azurewebsites.net/",
  "https://powerbi.microsoft.com/es-es/",
  "https://www.dropbox.com/",
  "https://dataintelligence.store/",
  "https://app-data-i.users.earthengine.app/",
  "https://odooutil.azurewebsites.net/",
  "https://filtradordashboard.azurewebsites.net/",
  "https://ide.dataintelligence-group.com/mapstore/#/",
  "https://ide.dataintelligence-group.com/geoserver/web/?0",
  "https://ide.dataintelligence-group.com/",
  "https://rpubs.com/dataintelligence/",
  "https://github.com/Sud-Austral/",
  "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$startRow = 338
for ($i = 0; $i -lt 14; $i++) {
  $row = $startRow + $i
  $name = $names[$i]
  $url = $urls[$i]

  $ws.Cells.Item($row, 1).Value = $name
  $ws.Cells.Item($row, 2).Value = $url
  $ws.Cells.Item($row, 3).Value = "Disponible"
  $ws.Cells.Item($row, 4).Value = $newTimestamp

  $hashIndex = $url.IndexOf("#")
  if ($hashIndex -ge 0) {
    $baseUrl = $url.Substring(0, $hashIndex)
    $subAddress = $url.Substring($hashIndex + 1)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $baseUrl, $subAddress)
  } else {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $url)
  }

  $ws.Cells.Item($row, 2).Style = "Hyperlink"
  $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
